# Applies the "Added one line to the report" commit:
#   1. Merge the three adjacent "Menlo" runs that spell out the
#      poetry-run command into a single run (same visible text).
#   2. Append a new "APA 7" list item after the "Lauren Gee's
#      mario-locate-objects code" reference, matching the existing
#      ListParagraph / numPr formatting.

$d = $word.ActiveDocument

# --- Change 1: collapse the split runs around the poetry-run command ---
$findText = [char]0x2018 + "poetry run nes_py --rom super-mario-bros.nes --mode human" + [char]0x2019 + " for human controlled Mario"
$d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $findText, 2)

# --- Change 2: add a new "APA 7" reference bullet after Lauren Gee's ---
$targetText = "Lauren Gee" + [char]0x2019 + "s mario-locate-objects code"
$targetIdx = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i++
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText) {
        $targetIdx = $i
    }
}

if ($targetIdx -gt 0) {
    $target = $d.Paragraphs.Item($targetIdx)
    $target.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($targetIdx + 1)
    $newPara.Range.InsertAfter("APA 7")
}
